$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.540118098258972
$ws.Range("B1").Value = 2.774966955184937
$ws.Range("C1").Value = 3.074009656906128
$ws.Range("D1").Value = 3.006079912185669
$ws.Range("E1").Value = 2.569772720336914
